$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.327636122703552
$ws.Range("B1").Value = 1.706321120262146
$ws.Range("C1").Value = 2.340502738952637
$ws.Range("D1").Value = 6.540104866027832
$ws.Range("E1").Value = 2.848648309707642
